# Update the marksheet's "correct / total" marks figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Row 11 ("Marking" row): Right count changes from 3 to 5
$ws.Range("B11").Value = 5

# Row 12 ("Total" row): Right total changes from 69 to 115,
# and the Max cell text changes from "64/84" to "115/140"
$ws.Range("B12").Value = 115
$ws.Range("E12").Value = "115/140"
